# Applies the changes described in the commit:
#  - remove the _FilterDatabase defined name (and worksheet AutoFilter)
#  - restyle the header row (white bold text on a solid/black fill, no explicit font name)
#  - turn off the table's banded row/column stripes
#  - narrow column C slightly
#  - update the sample data row values

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Listado Proyectos")

# --- Header row formatting (C3:F3) ---
$headerRange = $ws.Range("C3:F3")
$headerRange.Font.Name = "Calibri"
$headerRange.Font.Color = 16777215   # white
$headerRange.Interior.Color = 0      # black

# --- Table style: disable row/column stripes ---
$lo = $ws.ListObjects.Item(1)
$lo.ShowTableStyleRowStripes = $False
$lo.ShowTableStyleColumnStripes = $False

# --- Column C width: 23 -> 22 characters ---
$ws.Columns.Item(3).ColumnWidth = 21.17

# --- Update sample data row (row 4) ---
$ws.Range("C4").Value = "Proyecto prueba 2"
$ws.Range("D4").Value = "añsdfkjañsldasñdlfkj"
$ws.Range("E4").Value = 2000000000
$ws.Range("F4").Value = "21/01/2025"

# --- Remove the worksheet AutoFilter and the _FilterDatabase defined name ---
$ws.AutoFilterMode = $False
foreach ($n in $wb.Names) {
  $n.Delete()
}
